# Update the weekly fruit/vegetable price data (Jengibre sheet).
# Columns: D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
#          M=Precio promedio ponderado, P=Precio $/Kg
# Only the rows below change; rows 6 and 8 remain the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = 44469; J = 140; K = 13000; L = 14000; M = 13500; P = 1038 },
    @{ Row = 3; D = 44406; J = 160; K = 17000; L = 18000; M = 17500; P = 1346 },
    @{ Row = 4; D = 44159; J = 100; K = 23000; L = 24000; M = 23500; P = 1808 },
    @{ Row = 5; D = 44397; J = 140; K = 12500; L = 13000; M = 12750; P = 981 },
    @{ Row = 7; D = 44389; J = 120; K = 12000; L = 13000; M = 12500; P = 962 },
    @{ Row = 9; D = 44229; J = 120; K = 44000; L = 45000; M = 44500; P = 3423 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("J$r").Value = $u.J
    $ws.Range("K$r").Value = $u.K
    $ws.Range("L$r").Value = $u.L
    $ws.Range("M$r").Value = $u.M
    $ws.Range("P$r").Value = $u.P
}
